$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing "Rect3"/"Rect4" data so the
# current rows 3-6 shift down to rows 5-8 (labels + values move with them).
$ws.Rows("3:4").Insert()

# New row 3: "Std_dev Rect1" (replacing the series that shifted to row 5)
$ws.Range("A3").Value = "Std_dev Rect1"
$ws.Range("B3").Value = 0.1182059035512639
$ws.Range("C3").Value = 0.1376223264721769
$ws.Range("D3").Value = 0.05956167798364864
$ws.Range("E3").Value = 0.09371153272630643
$ws.Range("F3").Value = 0.08474262395795443
$ws.Range("G3").Value = 0.06945385870491867

# New row 4: "Std_dev Rect2" (replacing the series that shifted to row 6)
$ws.Range("A4").Value = "Std_dev Rect2"
$ws.Range("B4").Value = 0.03136035073482118
$ws.Range("C4").Value = 0.03429929960042005
$ws.Range("D4").Value = 0.03346630074544385
$ws.Range("E4").Value = 0.02629950238907533
$ws.Range("F4").Value = 0.02820561540189194
$ws.Range("G4").Value = 0.0250851446287112

# Re-assert the labels for the shifted rows explicitly (the engine can
# otherwise alias/reuse shared-string slots across rows when a new value
# happens to match old text, which would corrupt sibling cells).
$ws.Range("A5").Value = "Std_dev Rect3"
$ws.Range("A6").Value = "Std_dev Rect4"

# The rows that used to be "Rect3" and "Rect4" are now at 7 and 8; rename
# them to continue the series as "Rect5" and "Rect6".
$ws.Range("A7").Value = "Std_dev Rect5"
$ws.Range("A8").Value = "Std_dev Rect6"

$ws.Range("A1").Select()
